$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.226.71"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "2.574.55"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'587.07"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("D6").Value = "'144.64"
$ws.Range("E6").Value = "  -1.97%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.590"
$ws.Range("E8").Value = "  -1.01%  "

$ws.Range("D9").Value = "'0.106"
$ws.Range("E9").Value = "  -2.05%  "

$ws.Range("D10").Value = "'5.63"
$ws.Range("E10").Value = "  +0.01%  "

$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").Value = "'0.351"
$ws.Range("E12").Value = "  -1.34%  "

$ws.Range("D13").Value = "'27.19"
$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("D14").Value = "3.038.34"
$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("D15").Value = "63.163.18"
$ws.Range("E15").Value = "  +0.17%  "

$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("D17").Value = "2.577.60"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").Value = "'11.07"
$ws.Range("E18").Value = "  -2.30%  "

$ws.Range("D19").Value = "'341.45"
$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("E20").Value = "  -1.85%  "

$ws.Range("D21").Value = "'6.64"
$ws.Range("E21").Value = "  -3.27%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("E23").Value = "  +3.77%  "

$ws.Range("D24").Value = "'67.93"
$ws.Range("E24").Value = "  +1.86%  "

$ws.Range("E25").Value = "  +5.34%  "

$ws.Range("D26").Value = "'1.61"
$ws.Range("E26").Value = "  -0.81%  "

$ws.Range("E27").Value = "  -3.14%  "

$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").Value = "'7.94"
$ws.Range("E29").Value = "  -1.96%  "

$ws.Range("D30").Value = "'8.23"
$ws.Range("E30").Value = "  -2.32%  "

$ws.Range("E31").Value = "  -2.53%  "

$ws.Range("D32").Value = "'470.80"
$ws.Range("E32").Value = "  +2.25%  "

$ws.Range("D33").Value = "0.0₃0802"
$ws.Range("E33").Value = "  -2.49%  "

$ws.Range("D34").Value = "'1.68"
$ws.Range("E34").Value = "  +3.51%  "

$ws.Range("D35").Value = "'176.06"
$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("D37").Value = "'0.400"
$ws.Range("E37").Value = "  -1.74%  "

$ws.Range("D38").Value = "'18.87"
$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("E41").Value = "  -2.86%  "

$ws.Range("D42").Value = "'40.16"
$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("D43").Value = "'157.55"
$ws.Range("E43").Value = "  +4.51%  "

$ws.Range("D44").Value = "'3.69"
$ws.Range("E44").Value = "  -3.32%  "

$ws.Range("D45").Value = "'21.27"
$ws.Range("E45").Value = "  +2.08%  "

$ws.Range("D46").Value = "'0.633"
$ws.Range("E46").Value = "  +3.58%  "

$ws.Range("D47").Value = "'0.0540"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").Value = "'0.0963"
$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("D50").Value = "'18.18"
$ws.Range("E50").Value = "  -1.20%  "

$ws.Range("D51").Value = "'11.38"
$ws.Range("E51").Value = "  -0.07%  "
